$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '43.168.11' },
    @{ Cell = 'E2'; Value = '  +4.44%  ' },
    @{ Cell = 'D3'; Value = '2.250.06' },
    @{ Cell = 'E3'; Value = '  +3.33%  ' },
    @{ Cell = 'E4'; Value = '  +0.13%  ' },
    @{ Cell = 'D5'; Value = '244.83' },
    @{ Cell = 'E5'; Value = '  +2.85%  ' },
    @{ Cell = 'E6'; Value = '  +1.42%  ' },
    @{ Cell = 'D7'; Value = '76.00' },
    @{ Cell = 'E7'; Value = '  +8.76%  ' },
    @{ Cell = 'E8'; Value = '  -0.06%  ' },
    @{ Cell = 'E9'; Value = '  +6.57%  ' },
    @{ Cell = 'D10'; Value = '41.31' },
    @{ Cell = 'E10'; Value = '  +4.81%  ' },
    @{ Cell = 'D11'; Value = '0.0938' },
    @{ Cell = 'D12'; Value = '7.04' },
    @{ Cell = 'E12'; Value = '  +4.81%  ' },
    @{ Cell = 'E13'; Value = '  +0.80%  ' },
    @{ Cell = 'D14'; Value = '2.592.08' },
    @{ Cell = 'E14'; Value = '  +3.58%  ' },
    @{ Cell = 'D15'; Value = '14.68' },
    @{ Cell = 'E15'; Value = '  +3.64%  ' },
    @{ Cell = 'D16'; Value = '2.261.62' },
    @{ Cell = 'E16'; Value = '  +4.30%  ' },
    @{ Cell = 'D17'; Value = '0.805' },
    @{ Cell = 'E17'; Value = '  +1.46%  ' },
    @{ Cell = 'D18'; Value = '43.096.93' },
    @{ Cell = 'E18'; Value = '  +4.79%  ' },
    @{ Cell = 'E19'; Value = '  +4.74%  ' },
    @{ Cell = 'D20'; Value = '71.48' },
    @{ Cell = 'E20'; Value = '  +1.29%  ' },
    @{ Cell = 'D21'; Value = '6.00' },
    @{ Cell = 'E21'; Value = '  +2.26%  ' },
    @{ Cell = 'D22'; Value = '9.99' },
    @{ Cell = 'E22'; Value = '  +6.09%  ' },
    @{ Cell = 'D23'; Value = '230.74' },
    @{ Cell = 'E23'; Value = '  +2.35%  ' },
    @{ Cell = 'D24'; Value = '2.19' },
    @{ Cell = 'E24'; Value = '  +15.55%  ' },
    @{ Cell = 'E25'; Value = '  -0.02%  ' },
    @{ Cell = 'D26'; Value = '10.96' },
    @{ Cell = 'E26'; Value = '  +2.05%  ' },
    @{ Cell = 'D27'; Value = '3.50' },
    @{ Cell = 'E27'; Value = '  +0.67%  ' },
    @{ Cell = 'B28'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D28'; Value = '39.25' },
    @{ Cell = 'E28'; Value = '  +30.52%  ' },
    @{ Cell = 'B29'; Value = 'PancakeSwap' },
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' },
    @{ Cell = 'D29'; Value = '2.25' },
    @{ Cell = 'E29'; Value = '  +2.65%  ' },
    @{ Cell = 'B30'; Value = 'Toncoin' },
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Cell = 'D30'; Value = '2.22' },
    @{ Cell = 'E30'; Value = '  +2.29%  ' },
    @{ Cell = 'D31'; Value = '173.51' },
    @{ Cell = 'E31'; Value = '  +3.56%  ' },
    @{ Cell = 'D32'; Value = '20.38' },
    @{ Cell = 'E32'; Value = '  +2.36%  ' },
    @{ Cell = 'D33'; Value = '0.0799' },
    @{ Cell = 'E33'; Value = '  +5.12%  ' },
    @{ Cell = 'E34'; Value = '  +5.24%  ' },
    @{ Cell = 'E35'; Value = '  +1.89%  ' },
    @{ Cell = 'E36'; Value = '  +10.95%  ' },
    @{ Cell = 'E37'; Value = '  +6.86%  ' },
    @{ Cell = 'D38'; Value = '0.0335' },
    @{ Cell = 'E38'; Value = '  +18.84%  ' },
    @{ Cell = 'D39'; Value = '13.33' },
    @{ Cell = 'E39'; Value = '  +13.37%  ' },
    @{ Cell = 'D40'; Value = '2.14' },
    @{ Cell = 'E40'; Value = '  +3.80%  ' },
    @{ Cell = 'D41'; Value = '5.55' },
    @{ Cell = 'E41'; Value = '  +3.28%  ' },
    @{ Cell = 'E42'; Value = '  +8.15%  ' },
    @{ Cell = 'D43'; Value = '60.18' },
    @{ Cell = 'E43'; Value = '  +2.20%  ' },
    @{ Cell = 'D44'; Value = '106.10' },
    @{ Cell = 'E44'; Value = '  +9.02%  ' },
    @{ Cell = 'D45'; Value = '8.74' },
    @{ Cell = 'E45'; Value = '  +5.85%  ' },
    @{ Cell = 'D46'; Value = '0.489' },
    @{ Cell = 'E47'; Value = '  +2.86%  ' },
    @{ Cell = 'E48'; Value = '  +10.53%  ' },
    @{ Cell = 'E49'; Value = '  +3.66%  ' },
    @{ Cell = 'E50'; Value = '  +3.30%  ' },
    @{ Cell = 'D51'; Value = '2.465.31' },
    @{ Cell = 'E51'; Value = '  +3.83%  ' }
)

foreach ($change in $changes) {
    $ref = $change.Cell
    $newValue = $change.Value
    $rowNum = $ref -replace "[A-Z]+", ""
    $styleSource = "B" + $rowNum
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $newValue
    $ws.Range($ref).Style = $ws.Range($styleSource).Style
}

Write-Output ("Applied " + $changes.Count + " cell updates")